$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived data for rows 2-7 (Sending cluster is now only FAPs/MuSCs;
# "ECs" sending-cluster rows were dropped and all numeric columns were
# recomputed with the new TPM values).
$rowsData = @(
  @("FAPs",  "Rspo1", "Znrf3", "ECs",   2, 0.6666666666666666, 1.52558,  4.57674,   0.8891290081558957, 0.8891290081558957, 3, 1, 0.5762313333333333, 1.728694,   0.07066599920647125, 0.07066599920647126, 0.8790869975066666, 7.91178297756,   0.0628311897847951,  0.0628311897847951),
  @("FAPs",  "Rspo1", "Znrf3", "FAPs",  2, 0.6666666666666666, 1.52558,  4.57674,   0.8891290081558957, 0.8891290081558957, 3, 1, 3.933615666666666, 11.800847,  0.4823980674067757,  0.4823980674067757,  6.001045388753333,  54.00940849878,  0.4289141152097073,  0.4289141152097074),
  @("FAPs",  "Rspo1", "Znrf3", "MuSCs", 2, 0.6666666666666666, 1.52558,  4.57674,   0.8891290081558957, 0.8891290081558957, 3, 1, 3.644447,           10.933341,  0.446935933386753,   0.446935933386753,   5.559895454259999,  50.03905908834,  0.3973837031613932,  0.3973837031613932),
  @("MuSCs", "Rspo1", "Znrf3", "ECs",   3, 1,                  0.190234, 0.570702,  0.1108709918441043, 0.1108709918441043, 3, 1, 0.5762313333333333, 1.728694,   0.07066599920647125, 0.07066599920647126, 0.1096187914653333, 0.9865691231880001, 0.007834809421676157, 0.007834809421676159),
  @("MuSCs", "Rspo1", "Znrf3", "FAPs",  3, 1,                  0.190234, 0.570702,  0.1108709918441043, 0.1108709918441043, 3, 1, 3.933615666666666, 11.800847,  0.4823980674067757,  0.4823980674067757,  0.7483074427326667, 6.734766984594,  0.05348395219706831, 0.05348395219706832),
  @("MuSCs", "Rspo1", "Znrf3", "MuSCs", 3, 1,                  0.190234, 0.570702,  0.1108709918441043, 0.1108709918441043, 3, 1, 3.644447,           10.933341,  0.446935933386753,   0.446935933386753,   0.6932977305979999, 6.239679575382,  0.04955223022535984, 0.04955223022535985)
)

for ($r = 0; $r -lt $rowsData.Length; $r++) {
  $rowArr = $rowsData[$r]
  $excelRow = $r + 2
  for ($c = 0; $c -lt $rowArr.Length; $c++) {
    $ws.Cells.Item($excelRow, $c + 1).Value = $rowArr[$c]
  }
}

# The refreshed data only has 6 data rows (was 9): drop the now-unused rows 8-10.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
